# Update loading_percent values for rows 2-25 (Case 5_92, 380 kV)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 7.965722762100372
    "C2" = 5.141302184792995
    "D2" = 5.033115134359028
    "E2" = 12.80736667787291
    "F2" = 24.52540140325567
    "I2" = 20.88698452704859
    "K2" = 8.019165173221673
    "M2" = 13.34531255437553
    "O2" = 21.98715207275361
    "B3" = 7.65402397521884
    "C3" = 4.955890905370566
    "D3" = 4.977371347440099
    "E3" = 12.60303087404094
    "F3" = 24.55218966510245
    "I3" = 20.97441576261756
    "K3" = 7.732999789525079
    "M3" = 13.16692806278943
    "O3" = 22.05587897094485
    "B4" = 7.456945490711584
    "C4" = 4.837442103993507
    "D4" = 4.942317558238714
    "E4" = 12.48047349263391
    "F4" = 24.57542221950041
    "I4" = 21.03228952619605
    "K4" = 7.550151859790222
    "M4" = 13.05933914864439
    "O4" = 22.10305660923367
    "B5" = 7.375323452522399
    "C5" = 4.788058318118646
    "D5" = 4.927831937099719
    "E5" = 12.43133169749129
    "F5" = 24.58659232514357
    "I5" = 21.05692594679661
    "K5" = 7.476200280888029
    "M5" = 13.01603650315331
    "O5" = 22.12353069553708
    "B6" = 7.361694797733451
    "C6" = 4.779792200726621
    "D6" = 4.925414683452067
    "E6" = 12.4232221718813
    "F6" = 24.58854983368147
    "I6" = 21.06108032764654
    "K6" = 7.467170679401572
    "M6" = 13.00888031663801
    "O6" = 22.12700573098022
    "B7" = 7.455849845605078
    "C7" = 4.836780550537584
    "D7" = 4.942123004240051
    "E7" = 12.47980741243562
    "F7" = 24.57556597457767
    "I7" = 21.03261752230942
    "K7" = 7.549130545273497
    "M7" = 13.05875289551951
    "O7" = 22.10332767787493
    "B8" = 7.859499201187104
    "C8" = 5.078352582384667
    "D8" = 5.014070449092727
    "E8" = 12.7363514881349
    "F8" = 24.53322811958208
    "I8" = 20.91626022354453
    "K8" = 7.922014561031979
    "M8" = 13.28343496710809
    "O8" = 22.00981381961487
    "B9" = 8.601234424360706
    "C9" = 5.513886891621224
    "D9" = 5.148299990107509
    "E9" = 13.25932287702087
    "F9" = 24.5041484902641
    "I9" = 20.72139623791939
    "K9" = 8.594070039828349
    "M9" = 13.73706833516162
    "O9" = 21.8660781323242
    "B10" = 9.110477561380433
    "C10" = 5.808712712280439
    "D10" = 5.242351914854596
    "E10" = 13.65140948460381
    "F10" = 24.51576988589884
    "I10" = 20.59862204425916
    "K10" = 9.048888468086641
    "M10" = 14.07512027044102
    "O10" = 21.7848251991386
    "B11" = 9.333490420013463
    "C11" = 5.937057401422086
    "D11" = 5.284069260736397
    "E11" = 13.83059310120407
    "F11" = 24.52822101173201
    "I11" = 20.54721539386344
    "K11" = 9.246858962109755
    "M11" = 14.2292756108704
    "O11" = 21.75318400529708
    "B12" = 9.41663473316658
    "C12" = 5.984806324020411
    "D12" = 5.299706755940965
    "E12" = 13.89849095896478
    "F12" = 24.53396456655485
    "I12" = 20.52838970739786
    "K12" = 9.320508049894711
    "M12" = 14.28764930269591
    "O12" = 21.74197012221235
    "B13" = 9.398787100625398
    "C13" = 5.974560972415569
    "D13" = 5.296346156374822
    "E13" = 13.88386708427772
    "F13" = 24.53268188013465
    "I13" = 20.53241562197838
    "K13" = 9.304705513978742
    "M13" = 14.27507842103518
    "O13" = 21.74435103576284
    "B14" = 9.340357251516409
    "C14" = 5.94100293431017
    "D14" = 5.2853590089456
    "E14" = 13.83617860997131
    "F14" = 24.52867293627362
    "I14" = 20.54565373807707
    "K14" = 9.252944710849329
    "M14" = 14.23407839364857
    "O14" = 21.75224602635061
    "B15" = 9.304395513631704
    "C15" = 5.920336048848774
    "D15" = 5.27860804438226
    "E15" = 13.80697175704954
    "F15" = 24.52635122747258
    "I15" = 20.55384598555022
    "K15" = 9.221067145753885
    "M15" = 14.20896286649522
    "O15" = 21.75718202471966
    "B16" = 9.095723630263318
    "C16" = 5.80020694011136
    "D16" = 5.239603551682131
    "E16" = 13.63970963483612
    "F16" = 24.51510022490802
    "I16" = 20.60207117220595
    "K16" = 9.035767936540315
    "M16" = 14.06504861764753
    "O16" = 21.78700034265922
    "B17" = 8.965448629305053
    "C17" = 5.725016138413486
    "D17" = 5.215397506003725
    "E17" = 13.53725736105801
    "F17" = 24.51003224200009
    "I17" = 20.63279528418393
    "K17" = 8.919780200733653
    "M17" = 13.97681784869508
    "O17" = 21.80665782983479
    "B18" = 8.889707392374175
    "C18" = 5.68122610695287
    "D18" = 5.201374738014733
    "E18" = 13.47841109920339
    "F18" = 24.50779179285801
    "I18" = 20.65088515368041
    "K18" = 8.852228425000385
    "M18" = 13.92610826829916
    "O18" = 21.8184649968999
    "B19" = 8.863925565900839
    "C19" = 5.666307150132949
    "D19" = 5.196609880101021
    "E19" = 13.45850304328205
    "F19" = 24.50714909866641
    "I19" = 20.65708183776948
    "K19" = 8.829213626259486
    "M19" = 13.90894718262751
    "O19" = 21.82254860889828
    "B20" = 8.979400992261066
    "C20" = 5.733076621983919
    "D20" = 5.217984687531703
    "E20" = 13.54815567507464
    "F20" = 24.51050193583573
    "I20" = 20.62948135721961
    "K20" = 8.932214357619381
    "M20" = 13.98620659300736
    "O20" = 21.804513415427
    "B21" = 9.357555410732243
    "C21" = 5.950883040073342
    "D21" = 5.288590592153461
    "E21" = 13.85018522412683
    "F21" = 24.52982256314449
    "I21" = 20.54174797383127
    "K21" = 9.268184132916788
    "M21" = 14.24612157299842
    "O21" = 21.74990621318275
    "B22" = 9.597066772581815
    "C22" = 6.088255351276457
    "D22" = 5.333800434807985
    "E22" = 14.04780851368947
    "F22" = 24.54844358428097
    "I22" = 20.48814564366136
    "K22" = 9.480064684582439
    "M22" = 14.41595712097
    "O22" = 21.71869459519715
    "B23" = 9.469952109593843
    "C23" = 6.015399099753391
    "D23" = 5.309758740456396
    "E23" = 13.94233582018649
    "F23" = 24.53795757022235
    "I23" = 20.51641170678648
    "K23" = 9.367694209941277
    "M23" = 14.32533371420371
    "O23" = 21.73494228190328
    "B24" = 8.973095754246863
    "C24" = 5.729434225674114
    "D24" = 5.216815353643929
    "E24" = 13.54322837066202
    "F24" = 24.51028749014063
    "I24" = 20.63097825622492
    "K24" = 8.926595579064319
    "M24" = 13.98196188929637
    "O24" = 21.80548132977871
    "B25" = 8.406505697084194
    "C25" = 5.400362215860368
    "D25" = 5.11276213820237
    "E25" = 13.11615760077365
    "F25" = 24.50622602851737
    "I25" = 20.77053715954404
    "K25" = 8.418916833089558
    "M25" = 13.61328190770182
    "O25" = 21.9006988495608
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
